$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (row, new text)
$updates = @(
    @{Row = 1;  Text = "0M"},
    @{Row = 2;  Text = "0M"},
    @{Row = 3;  Text = "0M"},
    @{Row = 4;  Text = "33"},
    @{Row = 7;  Text = "0.00018"},
    @{Row = 8;  Text = "0.00006"},
    @{Row = 9;  Text = "0.00026"},
    @{Row = 10; Text = "0.00034"},
    @{Row = 11; Text = "0.00048"},
    @{Row = 12; Text = "0.00747"},
    @{Row = 44; Text = "99.99"},
    @{Row = 45; Text = "0.01"},
    @{Row = 46; Text = "65"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, 1)
    $cell.Range.Text = $u.Text
}
